# Add the Week 12 update text as three separate runs (matching how Word
# would have recorded three distinct typing/edit actions), right after the
# existing "Week 12 :" run, inside the same paragraph.

$d = $word.ActiveDocument

# 1) Find the paragraph that currently contains only "Week 12 :"
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    $t = $t.TrimEnd([char]13, [char]7)   # strip trailing pilcrow / cell mark
    if ($t -eq "Week 12 :") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not find the 'Week 12 :' paragraph"
}

# 2) Pull that paragraph's real OOXML (<w:p ...>...</w:p>) out of the
#    document's full WordOpenXML so we keep its existing paraId / rsid /
#    pPr attributes untouched.
$full = $d.Content.WordOpenXML
$paraRe = [regex]'(?s)<w:p\b.*?</w:p>'
$paraXml = $null
foreach ($m in $paraRe.Matches($full)) {
    $textOnly = [regex]::Replace($m.Value, '<[^>]+>', '')
    if ($textOnly -eq "Week 12 :") {
        $paraXml = $m.Value
        break
    }
}
if ($paraXml -eq $null) {
    throw "Could not isolate the 'Week 12 :' paragraph XML"
}

# 3) Build the three new runs, each with the same run formatting as the
#    existing "Week 12 :" run (Calibri, sz 24 / szCs 24).
$newRuns = @(
    ' Icon & Splashscreen veranderd',
    ', navigatie orderlijker gemaakt en',
    ' UI aanpassingen (3u)'
)

$rPr = '<w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

$runsXml = ""
foreach ($t in $newRuns) {
    $escaped = $t.Replace("&", "&amp;")
    if ($t.StartsWith(" ") -or $t.EndsWith(" ")) {
        $spaceAttr = ' xml:space="preserve"'
    } else {
        $spaceAttr = ""
    }
    $runsXml += "<w:r>$rPr<w:t$spaceAttr>$escaped</w:t></w:r>"
}

# 4) Splice the new runs in just before the paragraph's closing tag.
$newParaXml = $paraXml -replace '</w:p>$', ($runsXml + '</w:p>')

# 5) Feed the modified paragraph back in via InsertXML, replacing the whole
#    paragraph range so the separate <w:r> elements survive (a plain
#    InsertAfter would normalize/merge runs that end up with identical
#    formatting).
$xmlDoc = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
$newParaXml
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$target.Range.InsertXML($xmlDoc)
